$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert rows so PIPO5/6/7 each have 3 replicate rows like the rest ---
# Insert 2 rows after row 14 (PIPO5 block) -> new rows become 15 and 16
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(15).Insert()

# Insert 2 rows after row 17 (PIPO6, originally row 15) -> new rows become 18 and 19
$ws.Rows.Item(18).Insert()
$ws.Rows.Item(18).Insert()

# Insert 2 rows after row 20 (PIPO7, originally row 16) -> new rows become 21 and 22
$ws.Rows.Item(21).Insert()
$ws.Rows.Item(21).Insert()

# --- Fill in tree IDs for the newly inserted rows ---
$ws.Range("A15").Value = $ws.Range("A14").Value2
$ws.Range("A16").Value = $ws.Range("A14").Value2
$ws.Range("A18").Value = $ws.Range("A17").Value2
$ws.Range("A19").Value = $ws.Range("A17").Value2
$ws.Range("A21").Value = $ws.Range("A20").Value2
$ws.Range("A22").Value = $ws.Range("A20").Value2

# --- Fill in measurement data (stomatal density raw counts) ---
$ws.Range("B14").Value = 1.657
$ws.Range("C14").Value = 79
$ws.Range("B15").Value = 0.577
$ws.Range("C15").Value = 27
$ws.Range("B16").Value = 0.495
$ws.Range("C16").Value = 25

$ws.Range("B17").Value = 2.463
$ws.Range("C17").Value = 105
$ws.Range("B18").Value = 0.875
$ws.Range("C18").Value = 42
$ws.Range("B19").Value = 0.404
$ws.Range("C19").Value = 23

$ws.Range("B20").Value = 1.727
$ws.Range("C20").Value = 71
$ws.Range("B21").Value = 1.106
$ws.Range("C21").Value = 43
$ws.Range("B22").Value = 0.7
$ws.Range("C22").Value = 32

# --- Formulas for the D column (C/B ratio) in the newly inserted rows ---
$ws.Range("D15").Formula = "=C15/B15"
$ws.Range("D16").Formula = "=C16/B16"
$ws.Range("D18").Formula = "=C18/B18"
$ws.Range("D19").Formula = "=C19/B19"
$ws.Range("D21").Formula = "=C21/B21"
$ws.Range("D22").Formula = "=C22/B22"

# --- Bold the TreeID column for all the trees that now have complete (3-rep) data ---
$ws.Range("A2:A22").Font.Bold = $true

Write-Output "done filling data"
